$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update duration values (column C) for the changed tasks
$ws.Range("C5").Value = 30    # win conditions: 40 -> 30
$ws.Range("C7").Value = 15    # comment all: 20 -> 15
$ws.Range("C9").Value = 0     # analysis: 90 -> 0
$ws.Range("C11").Value = 40   # test plan: 45 -> 40
$ws.Range("C12").Value = 115  # testing: 120 -> 115

# Recalculate formulas (SUM and hours conversion) so cached values refresh
$excel.Calculate()

# Update the active selection on the sheet
$ws.Range("G10").Select()
